$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1).Range
$end = $p1.End - 1

$r1 = $d.Range($end, $end)
$r1.InsertAfter(" (")

$end2 = $end + 2
$r2 = $d.Range($end2, $end2)
$r2.InsertAfter("Changed main")

$end3 = $end2 + 12
$r3 = $d.Range($end3, $end3)
$r3.InsertAfter(")")
